$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric are stored as literal text,
# matching the source data which keeps values like "18.00" / "1.000" as strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.240.20'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.909.09'
$ws.Range("E3").Value = '  +1.87%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.70'
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5068'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3930'
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("E9").Value = '  -4.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.141'
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.84'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.409'
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.915.67'
$ws.Range("E13").Value = '  +2.32%  '
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.89'
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.312'
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.71'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06608'
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.00'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.203'
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.299.45'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").Value = '  +0.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.317'
$ws.Range("E25").Value = '  +1.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.596'
$ws.Range("E26").Value = '  +1.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.133.46'
$ws.Range("E27").Value = '  +2.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.05'
$ws.Range("E28").Value = '  -0.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '158.02'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.35'
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.103'
$ws.Range("E31").Value = '  +3.02%  '
$ws.Range("E32").Value = '  +0.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.646'
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.615'
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.682'
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06660'
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02415'
$ws.Range("E37").Value = '  +1.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.248'
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2188'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.280'
$ws.Range("E40").Value = '  +7.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6429'
$ws.Range("E41").Value = '  +0.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.015'
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.30'
$ws.Range("E45").Value = '  -1.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6019'
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.721'
$ws.Range("E47").Value = '  +1.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.277'
$ws.Range("E48").Value = '  +1.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.021'
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '123.07'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.187'
$ws.Range("E51").Value = '  -0.88%  '
